$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab/title) from "BrassA-HW25.xpc" to "BrassA"
$ws.Name = "BrassA"

# Append a new data row (row 16) with the averaged intensities for the
# "HexGrid-60degTilt5degRes" scheme (HKL index 14)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.085420576528351
$ws.Range("D16").Value = 0.9436809254642285
$ws.Range("E16").Value = 0.986338913977008
$ws.Range("F16").Value = 0.975497084934187
$ws.Range("G16").Value = 1.085420576528351
$ws.Range("H16").Value = 0.9436809254642285
$ws.Range("I16").Value = 1.013507860391305
$ws.Range("J16").Value = 0.962453540013901
$ws.Range("K16").Value = 1.023791702434834
$ws.Range("L16").Value = 0.9516332046643771
$ws.Range("M16").Value = 1.085420576528351
$ws.Range("N16").Value = 0.9650099197206183
$ws.Range("O16").Value = 0.9977343752259437
$ws.Range("P16").Value = 0.992790476051024

# Match the formatting of the first (HKL index) column used throughout the
# table by copying the style from the cell directly above it.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
